$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order/values for the CHIMANIMANI institution rows (rows 2-8).
# Each institution keeps its own Masculino/Feminino/TOTAL counts; only the
# row order changes (institutions re-sorted), which shifts which counts
# line up with which row.
$rows = @(
    @{ Row = 2;  Inst = "UNIZAMBEZE";            M = 0;  F = 1;  T = 1  },
    @{ Row = 3;  Inst = "MICAIA";                 M = 3;  F = 3;  T = 6  },
    @{ Row = 4;  Inst = "SDAE SUSSUNDENGA";        M = 1;  F = 0;  T = 1  },
    @{ Row = 5;  Inst = "PARQUE DE CHIMANIMANI";   M = 1;  F = 0;  T = 1  },
    @{ Row = 6;  Inst = "ITAM";                   M = 1;  F = 0;  T = 1  },
    @{ Row = 7;  Inst = "ISPM";                   M = 0;  F = 1;  T = 1  },
    @{ Row = 8;  Inst = "UCM";                    M = 1;  F = 0;  T = 1  }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Inst
    $ws.Cells.Item($r.Row, 3).Value = $r.M
    $ws.Cells.Item($r.Row, 4).Value = $r.F
    $ws.Cells.Item($r.Row, 5).Value = $r.T
}
